# Update Name of Algo
# Applies updated RandomForest imputation result values to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C2"  = -12.1902
    "B3"  = 6.045200000000006
    "E3"  = 16.5369
    "E12" = 18.48930000000002
    "B14" = 5.960599999999999
    "B21" = 9.447699999999999
    "B23" = 9.311300000000003
    "E24" = 16.371
    "B25" = 5.671400000000002
    "C25" = -11.7215
    "E25" = 16.98100000000001
    "B26" = 5.018700000000006
    "C27" = -12.82789999999999
    "B29" = 5.038500000000002
    "C31" = -13.4254
    "C39" = -12.72480000000001
    "C48" = -11.49099999999999
    "E50" = 16.2873
    "C51" = -11.6318
    "C52" = -11.5524
    "B53" = 6.079700000000003
    "E53" = 16.49730000000001
    "C55" = -13.57479999999999
    "C56" = -13.22759999999999
    "B57" = 4.930499999999995
    "C57" = -13.60639999999999
    "E57" = 16.6166
    "B59" = 4.851700000000001
    "E61" = 16.56650000000001
    "E63" = 18.57410000000002
    "B69" = 5.652599999999993
    "E70" = 18.64420000000002
    "C73" = -13.2509
    "B79" = 9.262600000000008
    "B83" = 5.574899999999998
    "E86" = 16.5448
    "C89" = -9.975800000000005
    "C90" = -12.6003
    "B91" = 5.017100000000002
    "C92" = -10.3813
    "B93" = 5.675699999999998
    "E98" = 15.8295
    "E100" = 16.3623
    "E102" = 16.7493
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
